# Auto-generated Excel COM-interop script
# Applies targeted cell updates across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 5657
$ws.Range("J3").Value = 5657
$ws.Range("L3").Value = 5657
$ws.Range("N3").Value = -5885

$ws.Range("H12").Value = 327.8
$ws.Range("I12").Value = 284.75
$ws.Range("K12").Value = 284.75
$ws.Range("M12").Value = -114.75

$ws.Range("H98").Value = 1456.9375
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H102").Value = 5657
$ws.Range("J102").Value = 5657
$ws.Range("L102").Value = 5657
$ws.Range("N102").Value = -12147

$ws.Range("H116").Value = 3383.3333
$ws.Range("I116").Value = 3157.1428
$ws.Range("J116").Value = 3700
$ws.Range("K116").Value = 3157.1428
$ws.Range("L116").Value = 3700
$ws.Range("M116").Value = 284.8571999999999
$ws.Range("N116").Value = -10584

$ws.Range("H122").Value = 1456.9375
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 10757763
$ws.Range("I132").Value = 12823351
$ws.Range("J132").Value = 16700.2
$ws.Range("K132").Value = 38470053
$ws.Range("L132").Value = 50100.60000000001
$ws.Range("M132").Value = -38467523
$ws.Range("N132").Value = -55160.60000000001

$ws.Range("H135").Value = 186.33333
$ws.Range("I135").Value = 79.5
$ws.Range("K135").Value = 715.5
$ws.Range("M135").Value = 1819.5

$ws.Range("H138").Value = 2104.1313
$ws.Range("I138").Value = 1063
$ws.Range("J138").Value = 2247.7356
$ws.Range("K138").Value = 3189
$ws.Range("L138").Value = 6743.2068
$ws.Range("M138").Value = 1951
$ws.Range("N138").Value = -17023.2068

$ws.Range("H141").Value = 9018.429
$ws.Range("I141").Value = 9904.083000000001
$ws.Range("K141").Value = 29712.249
$ws.Range("M141").Value = -24532.249

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 644.8723
$ws.Range("I74").Value = 591.2632
$ws.Range("K74").Value = 591.2632
$ws.Range("M74").Value = 282.7368

$ws.Range("H77").Value = 644.8723
$ws.Range("I77").Value = 591.2632
$ws.Range("K77").Value = 2956.316
$ws.Range("M77").Value = 1411.684

$ws.Range("H97").Value = 420.9
$ws.Range("I97").Value = 420.9
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 420.9
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = 75.10000000000002
$ws.Range("N97").ClearContents()

$ws.Range("H103").Value = 60400
$ws.Range("J103").Value = 60400
$ws.Range("L103").Value = 60400
$ws.Range("N103").Value = -62744

$ws.Range("H132").Value = 4447
$ws.Range("I132").Value = 5030.75
$ws.Range("J132").Value = 3668.6667
$ws.Range("K132").Value = 15092.25
$ws.Range("L132").Value = 11006.0001
$ws.Range("M132").Value = -12562.25
$ws.Range("N132").Value = -16066.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H41").Value = 249266.67
$ws.Range("J41").Value = 249266.67
$ws.Range("L41").Value = 249266.67
$ws.Range("N41").Value = -250042.67

$ws.Range("H48").Value = 249266.67
$ws.Range("J48").Value = 249266.67
$ws.Range("L48").Value = 249266.67
$ws.Range("N48").Value = -250096.67

$ws.Range("H64").Value = 279.8
$ws.Range("I64").Value = 333
$ws.Range("K64").Value = 333
$ws.Range("M64").Value = -108

$ws.Range("H67").Value = 279.8
$ws.Range("I67").Value = 333
$ws.Range("K67").Value = 333
$ws.Range("M67").Value = 447

$ws.Range("H86").Value = 3508.818
$ws.Range("I86").Value = 3352.0952
$ws.Range("K86").Value = 3352.0952
$ws.Range("M86").Value = -2229.0952

$ws.Range("H89").Value = 3508.818
$ws.Range("I89").Value = 3352.0952
$ws.Range("K89").Value = 16760.476
$ws.Range("M89").Value = -11144.476

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 402
$ws.Range("I22").Value = 360
$ws.Range("J22").Value = 465
$ws.Range("K22").Value = 360
$ws.Range("L22").Value = 465
$ws.Range("M22").Value = -10
$ws.Range("N22").Value = -1165

$ws.Range("H31").Value = 816.6323
$ws.Range("I31").Value = 722.04346
$ws.Range("J31").Value = 1014.4091
$ws.Range("K31").Value = 722.04346
$ws.Range("L31").Value = 1014.4091
$ws.Range("M31").Value = -427.04346
$ws.Range("N31").Value = -1604.4091

$ws.Range("H34").Value = 816.6323
$ws.Range("I34").Value = 722.04346
$ws.Range("J34").Value = 1014.4091
$ws.Range("K34").Value = 722.04346
$ws.Range("L34").Value = 1014.4091
$ws.Range("M34").Value = -520.04346
$ws.Range("N34").Value = -1418.4091

$ws.Range("H132").Value = 8409.764999999999
$ws.Range("I132").Value = 13825.5
$ws.Range("K132").Value = 41476.5
$ws.Range("M132").Value = -38946.5

$ws.Range("H134").Value = 10102238
$ws.Range("I134").Value = 10753821
$ws.Range("K134").Value = 32261463
$ws.Range("M134").Value = -32258928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 3014.2856
$ws.Range("J22").Value = 3083.3333
$ws.Range("L22").Value = 9249.999899999999
$ws.Range("N22").Value = -9587.999899999999

$ws.Range("H27").Value = 3014.2856
$ws.Range("J27").Value = 3083.3333
$ws.Range("L27").Value = 9249.999899999999
$ws.Range("N27").Value = -9453.999899999999

$ws.Range("H68").Value = 2544
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2544
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 7632
$ws.Range("N68").Value = -9254
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 2544
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2544
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 22896
$ws.Range("N71").Value = -31008
$ws.Range("M71").ClearContents()

$ws.Range("H106").Value = 3905.8
$ws.Range("J106").Value = 3905.8
$ws.Range("L106").Value = 11717.4
$ws.Range("N106").Value = -13609.4

$ws.Range("H112").Value = 76936150
$ws.Range("J112").Value = 83347170
$ws.Range("L112").Value = 250041510
$ws.Range("N112").Value = -250043726

$ws.Range("H122").Value = 688.6
$ws.Range("I122").Value = 302.66666
$ws.Range("K122").Value = 2723.99994
$ws.Range("M122").Value = -273.9999399999997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 28643
$ws.Range("J93").Value = 28643
$ws.Range("L93").Value = 28643
$ws.Range("N93").Value = -32387

$ws.Range("H97").Value = 726.3333
$ws.Range("I97").Value = 742.25
$ws.Range("J97").Value = 599
$ws.Range("K97").Value = 742.25
$ws.Range("L97").Value = 599
$ws.Range("M97").Value = -246.25
$ws.Range("N97").Value = -1591

$ws.Range("H132").Value = 2766.28
$ws.Range("I132").Value = 2511.3076
$ws.Range("J132").Value = 3042.5
$ws.Range("K132").Value = 7533.9228
$ws.Range("L132").Value = 9127.5
$ws.Range("M132").Value = -5003.9228
$ws.Range("N132").Value = -14187.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1010.2381
$ws.Range("I16").Value = 969.1053000000001
$ws.Range("J16").Value = 1401
$ws.Range("K16").Value = 969.1053000000001
$ws.Range("L16").Value = 1401
$ws.Range("M16").Value = -799.1053000000001
$ws.Range("N16").Value = -1741

$ws.Range("H46").Value = 1433.3334
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 1900
$ws.Range("K46").Value = 500
$ws.Range("L46").Value = 1900
$ws.Range("M46").Value = -312
$ws.Range("N46").Value = -2276

$ws.Range("H100").Value = 1418.1818
$ws.Range("J100").Value = 1846.6666
$ws.Range("L100").Value = 1846.6666
$ws.Range("N100").Value = -2928.6666

$ws.Range("H132").Value = 52269.15
$ws.Range("I132").Value = 1855.7142
$ws.Range("J132").Value = 169900.5
$ws.Range("K132").Value = 5567.142599999999
$ws.Range("L132").Value = 509701.5
$ws.Range("M132").Value = -3037.142599999999
$ws.Range("N132").Value = -514761.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2251.3044
$ws.Range("I132").Value = 1829.1
$ws.Range("J132").Value = 5066
$ws.Range("K132").Value = 5487.299999999999
$ws.Range("L132").Value = 15198
$ws.Range("M132").Value = -2957.299999999999
$ws.Range("N132").Value = -20258

$ws.Range("H136").Value = 1677.3077
$ws.Range("I136").Value = 1480.1
$ws.Range("K136").Value = 4440.299999999999
$ws.Range("M136").Value = -1890.299999999999
